$wb = $excel.ActiveWorkbook

# --- Sheet: LP1912 ---
$ws = $wb.Worksheets.Item('LP1912')
$ws.Cells.Item(2, 1).Value = 'Última actualización: 08:40:53'
$ws.Cells.Item(3, 1).Value = 'Total filas: 124'
$ws.Cells.Item(8, 1).Value = '04:44:46'
$ws.Cells.Item(8, 3).Value = '15_ABASTO'
$ws.Cells.Item(8, 4).Value = 2
$ws.Cells.Item(9, 1).Value = '03:52:04'
$ws.Cells.Item(9, 3).Value = '215A_EL PATO'
$ws.Cells.Item(9, 4).Value = 54
$ws.Cells.Item(10, 3).Value = '215_EL PELIGRO'
$ws.Cells.Item(38, 3).Value = '16_SANTA ANA'
$ws.Cells.Item(39, 3).Value = '17X38_ROMERO'
$ws.Cells.Item(45, 1).Value = '05:16:02'
$ws.Cells.Item(45, 3).Value = '17_ROMERO'
$ws.Cells.Item(45, 4).Value = 94
$ws.Cells.Item(46, 1).Value = '06:46:37'
$ws.Cells.Item(46, 3).Value = '215A_EL PATO'
$ws.Cells.Item(46, 4).Value = 4
$ws.Cells.Item(49, 1).Value = '06:53:56'
$ws.Cells.Item(49, 3).Value = '17_ROMERO'
$ws.Cells.Item(49, 4).Value = 1
$ws.Cells.Item(50, 1).Value = '06:46:37'
$ws.Cells.Item(50, 3).Value = '14_ABASTO'
$ws.Cells.Item(50, 4).Value = 8
$ws.Cells.Item(77, 1).Value = '06:18:01'
$ws.Cells.Item(77, 3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(77, 4).Value = 102
$ws.Cells.Item(78, 1).Value = '06:46:37'
$ws.Cells.Item(78, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(78, 4).Value = 74
$ws.Cells.Item(86, 1).Value = '08:10:38'
$ws.Cells.Item(86, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(86, 4).Value = 23
$ws.Cells.Item(87, 1).Value = '08:29:58'
$ws.Cells.Item(87, 3).Value = '215C_EL PATO'
$ws.Cells.Item(87, 4).Value = 4
$ws.Cells.Item(91, 1).Value = '08:40:53'
$ws.Cells.Item(91, 3).Value = '10_OLMOS'
$ws.Cells.Item(91, 4).Value = 7
$ws.Cells.Item(92, 1).Value = '08:40:53'
$ws.Cells.Item(92, 2).Value = '08:47'
$ws.Cells.Item(92, 3).Value = '215A_EL PATO'
$ws.Cells.Item(92, 4).Value = 7
$ws.Cells.Item(94, 1).Value = '08:10:38'
$ws.Cells.Item(94, 2).Value = '08:48'
$ws.Cells.Item(94, 3).Value = '10_OLMOS'
$ws.Cells.Item(94, 4).Value = 38
$ws.Cells.Item(95, 1).Value = '08:29:58'
$ws.Cells.Item(95, 2).Value = '08:50'
$ws.Cells.Item(95, 4).Value = 21
$ws.Cells.Item(96, 1).Value = '08:40:53'
$ws.Cells.Item(96, 2).Value = '08:51'
$ws.Cells.Item(96, 3).Value = '16_P MOR-SANTA ANA'
$ws.Cells.Item(96, 4).Value = 11
$ws.Cells.Item(97, 1).Value = '08:40:53'
$ws.Cells.Item(97, 2).Value = '08:59'
$ws.Cells.Item(97, 3).Value = '215B_EL PATO'
$ws.Cells.Item(97, 4).Value = 19
$ws.Cells.Item(98, 1).Value = '08:10:38'
$ws.Cells.Item(98, 2).Value = '09:00'
$ws.Cells.Item(98, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(98, 4).Value = 50
$ws.Cells.Item(99, 1).Value = '08:40:53'
$ws.Cells.Item(99, 2).Value = '09:01'
$ws.Cells.Item(99, 3).Value = '16_SANTA ANA'
$ws.Cells.Item(99, 4).Value = 21
$ws.Cells.Item(100, 1).Value = '07:38:30'
$ws.Cells.Item(100, 2).Value = '09:02'
$ws.Cells.Item(100, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(100, 4).Value = 84
$ws.Cells.Item(102, 1).Value = '08:40:53'
$ws.Cells.Item(102, 2).Value = '09:03'
$ws.Cells.Item(102, 3).Value = '17X38_ROMERO'
$ws.Cells.Item(102, 4).Value = 23
$ws.Cells.Item(103, 1).Value = '08:40:53'
$ws.Cells.Item(103, 2).Value = '09:07'
$ws.Cells.Item(103, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(103, 4).Value = 27
$ws.Cells.Item(104, 1).Value = '08:10:38'
$ws.Cells.Item(104, 2).Value = '09:10'
$ws.Cells.Item(104, 3).Value = '27_EL RETIRO'
$ws.Cells.Item(104, 4).Value = 60
$ws.Cells.Item(105, 1).Value = '07:50:33'
$ws.Cells.Item(105, 2).Value = '09:12'
$ws.Cells.Item(105, 4).Value = 82
$ws.Cells.Item(106, 1).Value = '08:40:53'
$ws.Cells.Item(106, 2).Value = '09:14'
$ws.Cells.Item(106, 3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(106, 4).Value = 34
$ws.Cells.Item(107, 1).Value = '07:38:30'
$ws.Cells.Item(107, 2).Value = '09:15'
$ws.Cells.Item(107, 3).Value = '27_EL RETIRO'
$ws.Cells.Item(107, 4).Value = 97
$ws.Cells.Item(108, 1).Value = '08:40:53'
$ws.Cells.Item(108, 2).Value = '09:16'
$ws.Cells.Item(108, 3).Value = '27_EL RETIRO'
$ws.Cells.Item(108, 4).Value = 36
$ws.Cells.Item(109, 1).Value = '08:40:53'
$ws.Cells.Item(109, 2).Value = '09:18'
$ws.Cells.Item(109, 3).Value = '215_EL PELIGRO'
$ws.Cells.Item(109, 4).Value = 38
$ws.Cells.Item(110, 1).Value = '07:50:33'
$ws.Cells.Item(110, 2).Value = '09:19'
$ws.Cells.Item(110, 3).Value = '215_EL PELIGRO'
$ws.Cells.Item(110, 4).Value = 89
$ws.Cells.Item(111, 1).Value = '08:29:58'
$ws.Cells.Item(111, 2).Value = '09:26'
$ws.Cells.Item(111, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(111, 4).Value = 57
$ws.Cells.Item(112, 1).Value = '08:40:53'
$ws.Cells.Item(112, 2).Value = '09:28'
$ws.Cells.Item(112, 3).Value = '10_OLMOS'
$ws.Cells.Item(112, 4).Value = 48
$ws.Cells.Item(113, 1).Value = '08:40:53'
$ws.Cells.Item(113, 2).Value = '09:29'
$ws.Cells.Item(113, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(113, 4).Value = 49
$ws.Cells.Item(114, 1).Value = '08:10:38'
$ws.Cells.Item(114, 2).Value = '09:29'
$ws.Cells.Item(114, 3).Value = '10_OLMOS'
$ws.Cells.Item(114, 4).Value = 79
$ws.Cells.Item(115, 2).Value = '09:33'
$ws.Cells.Item(115, 4).Value = 64
$ws.Cells.Item(116, 1).Value = '08:40:53'
$ws.Cells.Item(116, 2).Value = '09:34'
$ws.Cells.Item(116, 4).Value = 54
$ws.Cells.Item(117, 1).Value = '08:40:53'
$ws.Cells.Item(117, 2).Value = '09:41'
$ws.Cells.Item(117, 3).Value = '16_SANTA ANA'
$ws.Cells.Item(117, 4).Value = 61
$ws.Cells.Item(118, 1).Value = '08:40:53'
$ws.Cells.Item(118, 2).Value = '09:44'
$ws.Cells.Item(118, 3).Value = '14_ABASTO'
$ws.Cells.Item(118, 4).Value = 64
$ws.Cells.Item(119, 2).Value = '09:48'
$ws.Cells.Item(119, 3).Value = '15_ABASTO'
$ws.Cells.Item(119, 4).Value = 79
$ws.Cells.Item(120, 1).Value = '08:10:38'
$ws.Cells.Item(120, 2).Value = '09:49'
$ws.Cells.Item(120, 3).Value = '15_ABASTO'
$ws.Cells.Item(120, 4).Value = 99
$ws.Cells.Item(121, 2).Value = '09:50'
$ws.Cells.Item(121, 3).Value = '16_P MOR-SANTA ANA'
$ws.Cells.Item(121, 4).Value = 81
$ws.Cells.Item(122, 1).Value = '08:40:53'
$ws.Cells.Item(122, 2).Value = '09:51'
$ws.Cells.Item(122, 3).Value = '16_P MOR-SANTA ANA'
$ws.Cells.Item(122, 4).Value = 71
$ws.Cells.Item(123, 1).Value = '08:40:53'
$ws.Cells.Item(123, 2).Value = '09:56'
$ws.Cells.Item(123, 3).Value = '10_OLMOS'
$ws.Cells.Item(123, 4).Value = 76
$ws.Cells.Item(123, 5).Value = 'LP1912'
$ws.Cells.Item(124, 1).Value = '08:40:53'
$ws.Cells.Item(124, 2).Value = '10:03'
$ws.Cells.Item(124, 3).Value = '215C_EL PATO'
$ws.Cells.Item(124, 4).Value = 83
$ws.Cells.Item(124, 5).Value = 'LP1912'
$ws.Cells.Item(125, 1).Value = '08:40:53'
$ws.Cells.Item(125, 2).Value = '10:08'
$ws.Cells.Item(125, 3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(125, 4).Value = 88
$ws.Cells.Item(125, 5).Value = 'LP1912'
$ws.Cells.Item(126, 1).Value = '08:40:53'
$ws.Cells.Item(126, 2).Value = '10:18'
$ws.Cells.Item(126, 3).Value = '17_ROMERO'
$ws.Cells.Item(126, 4).Value = 98
$ws.Cells.Item(126, 5).Value = 'LP1912'
$ws.Cells.Item(127, 1).Value = '08:40:53'
$ws.Cells.Item(127, 2).Value = '10:20'
$ws.Cells.Item(127, 3).Value = '10_OLMOS'
$ws.Cells.Item(127, 4).Value = 100
$ws.Cells.Item(127, 5).Value = 'LP1912'
$ws.Cells.Item(128, 1).Value = '08:40:53'
$ws.Cells.Item(128, 2).Value = '10:32'
$ws.Cells.Item(128, 3).Value = '14_ABASTO'
$ws.Cells.Item(128, 4).Value = 112
$ws.Cells.Item(128, 5).Value = 'LP1912'
$ws.Cells.Item(129, 1).Value = '08:40:53'
$ws.Cells.Item(129, 2).Value = '10:34'
$ws.Cells.Item(129, 3).Value = '15_ABASTO'
$ws.Cells.Item(129, 4).Value = 114
$ws.Cells.Item(129, 5).Value = 'LP1912'

# --- Sheet: LP1912-215 ---
$ws = $wb.Worksheets.Item('LP1912-215')
$ws.Cells.Item(2, 1).Value = 'Última actualización: 08:40:53'
$ws.Cells.Item(20, 1).Value = '08:40:53'
$ws.Cells.Item(20, 4).Value = 7
$ws.Cells.Item(22, 1).Value = '08:40:53'
$ws.Cells.Item(22, 4).Value = 19
$ws.Cells.Item(23, 1).Value = '08:40:53'
$ws.Cells.Item(23, 4).Value = 38
$ws.Cells.Item(25, 1).Value = '08:40:53'
$ws.Cells.Item(25, 4).Value = 83

# --- Sheet: 6203-6173 ---
$ws = $wb.Worksheets.Item('6203-6173')
$ws.Cells.Item(2, 1).Value = 'Última actualización: 08:40:53'
$ws.Cells.Item(14, 1).Value = '08:40:53'
$ws.Cells.Item(14, 4).Value = 11
$ws.Cells.Item(16, 1).Value = '08:40:53'
$ws.Cells.Item(16, 4).Value = 75
$ws.Cells.Item(17, 1).Value = '08:40:53'
$ws.Cells.Item(17, 4).Value = 90
$ws.Cells.Item(18, 1).Value = '08:40:53'
$ws.Cells.Item(18, 4).Value = 101
